# Apply MN-specific calibration updates to the
# "Share of Cargo Dist Transported that is New This Year" workbook.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsPsgr  = $wb.Worksheets.Item("SoCDTtiNTY-psgr")
$wsFrgt  = $wb.Worksheets.Item("SoCDTtiNTY-frgt")

# Header row on both data sheets got taller (word-wrap needed more room).
$wsPsgr.Rows.Item(1).RowHeight = 45
$wsFrgt.Rows.Item(1).RowHeight = 45

# SoCDTtiNTY-psgr: LDVs (row 2) share recalibrated from 0.076 to 0.0755.
$wsPsgr.Range("B2:H2").Value = 0.0755

# SoCDTtiNTY-frgt: HDVs (row 3) shares recalibrated per-column.
$wsFrgt.Range("B3").Value = 0.0219
$wsFrgt.Range("C3").Value = 0.0219
$wsFrgt.Range("D3").Value = 0.081
$wsFrgt.Range("E3").Value = 0.0245
$wsFrgt.Range("F3").Value = 0.0219
$wsFrgt.Range("G3").Value = 0.0219
$wsFrgt.Range("H3").Value = 0.0219

# Update the saved selections / active sheet to match the author's session:
# - psgr sheet ends with A1:H7 selected
# - frgt sheet becomes the active tab with D10 selected
$wsPsgr.Activate()
$wsPsgr.Range("A1:H7").Select()

$wsFrgt.Activate()
$wsFrgt.Range("D10").Select()
